$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-08-04 Monday" "2025-08-05 Tuesday"

Replace-Text "58×78=" "54×64="
Replace-Text "58×81=" "30×66="
Replace-Text "81×25=" "70×89="
Replace-Text "29×70=" "21×79="
Replace-Text "85×23=" "70×65="
Replace-Text "79×99=" "49×16="
Replace-Text "69×63=" "89×51="
Replace-Text "90×13=" "79×50="
Replace-Text "32×31=" "56×50="
Replace-Text "44×70=" "39×17="
Replace-Text "37×33=" "66×43="
Replace-Text "93×19=" "46×92="
Replace-Text "39×20=" "44×30="
Replace-Text "42×46=" "25×79="
Replace-Text "90×14=" "59×52="
Replace-Text "92×51=" "17×97="
Replace-Text "97×23=" "11×76="
Replace-Text "55×85=" "76×67="
Replace-Text "47×60=" "48×53="
Replace-Text "82×68=" "49×88="
Replace-Text "21×95=" "56×36="
Replace-Text "74×12=" "81×88="
Replace-Text "18×90=" "15×52="
Replace-Text "62×47=" "36×25="
Replace-Text "28×86=" "22×15="
